# ZBP_08_pocet_aktivit.xlsx update:
#  - a new survey wave "2. 3. 2021" is added as the last date column on both
#    sheets ("data" -> column Z, "pocetR" -> column Y)
#  - the "aktualizace" (last-updated) date embedded in the two summary
#    title strings moves from 23. 2. 2021 to 9. 3. 2021

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": add column Z ("2. 3. 2021") with % values for rows 2..58
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Copy formatting (bold / centered / bordered) from the existing last
# header cell (Y1) onto the new header cell (Z1), then set its value.
$wsData.Range("Y1").Copy()
$wsData.Range("Z1").PasteSpecial(-4122)  # xlPasteFormats
$wsData.Range("Z1").Value = "2. 3. 2021"

$dataValues = @(0.4,0.31,0.29,0.5600000000000001,0.26,0.18,0.44,0.32,0.24,0.27,0.34,0.39,0.3,0.33,0.37,0.38,0.35,0.27,0.48,0.3,0.22,0.58,0.23,0.19,0.46,0.29,0.25,0.3,0.36,0.34,0.43,0.3,0.27,0.32,0.38,0.3,0.46,0.26,0.28,0.44,0.28,0.28,0.36,0.34,0.3,0.44,0.3,0.26,0.37,0.31,0.32,0.38,0.29,0.33,0.34,0.37,0.29)

for ($i = 0; $i -lt $dataValues.Count; $i++) {
    $row = $i + 2
    $wsData.Cells.Item($row, 26).Value = $dataValues[$i]   # column 26 = Z
}

# ---------------------------------------------------------------------
# Sheet "pocetR": add column Y ("2. 3. 2021") with counts for rows 2..20
# and an empty trailing cell on the "total" row 21
# ---------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

$wsPocet.Range("X1").Copy()
$wsPocet.Range("Y1").PasteSpecial(-4122)  # xlPasteFormats
$wsPocet.Range("Y1").Value = "2. 3. 2021"

$pocetValues = @(2130,512,776,842,624,721,579,370,702,1058,691,727,712,1042,1088,1107,490,252,281)

for ($i = 0; $i -lt $pocetValues.Count; $i++) {
    $row = $i + 2
    $wsPocet.Cells.Item($row, 25).Value = $pocetValues[$i]  # column 25 = Y
}

# Row 21 keeps the same empty-text pattern as the other trailing cells
# (B21:X21 are empty "t=s" placeholders) - copy an existing empty cell
# (content + formatting) across so Y21 matches the row's pattern exactly.
$wsPocet.Range("X21").Copy()
$wsPocet.Range("Y21").PasteSpecial(-4104)  # xlPasteAll

# ---------------------------------------------------------------------
# Update the "aktualizace" date embedded in the two (chart-title style)
# summary strings. Neither string is bound to a visible cell any more,
# but Cells.Replace still rewrites it wherever it lives in the workbook.
# ---------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("23. 2. 2021", "9. 3. 2021", -4142, 2, $false, $false, $true, $true)
}
